$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-08-24 Saturday" "2024-08-25 Sunday"

Replace-Text "47÷4=11, 3" "32÷5=6, 2"
Replace-Text "98÷8=12, 2" "31÷5=6, 1"
Replace-Text "15÷9=1, 6" "58÷5=11, 3"
Replace-Text "74÷5=14, 4" "23÷2=11, 1"
Replace-Text "47÷7=6, 5" "43÷4=10, 3"

Replace-Text "35÷8=4, 3" "56÷6=9, 2"
Replace-Text "93÷8=11, 5" "18÷4=4, 2"
Replace-Text "96÷3=32, 0" "32÷6=5, 2"
Replace-Text "18÷6=3, 0" "43÷4=10, 3"
Replace-Text "93÷6=15, 3" "35÷7=5, 0"

Replace-Text "92÷2=46, 0" "78÷5=15, 3"
Replace-Text "26÷7=3, 5" "63÷4=15, 3"
Replace-Text "10÷4=2, 2" "49÷4=12, 1"
Replace-Text "93÷2=46, 1" "50÷7=7, 1"
Replace-Text "52÷2=26, 0" "33÷2=16, 1"

Replace-Text "53÷5=10, 3" "93÷2=46, 1"
Replace-Text "43÷8=5, 3" "27÷2=13, 1"
Replace-Text "88÷5=17, 3" "58÷9=6, 4"
Replace-Text "43÷3=14, 1" "74÷4=18, 2"
Replace-Text "25÷3=8, 1" "38÷2=19, 0"

Replace-Text "36÷8=4, 4" "89÷2=44, 1"
Replace-Text "99÷7=14, 1" "96÷8=12, 0"
Replace-Text "35÷3=11, 2" "47÷3=15, 2"
Replace-Text "45÷5=9, 0" "59÷7=8, 3"
Replace-Text "45÷8=5, 5" "45÷6=7, 3"

Write-Output "Replacements complete"
